$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 111479727
$ws.Range("B4").Value = 56414
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("P4").Value = "SV Gärdvattnet, Jmt"
$ws.Range("Q4").Value = 506615.9431346679
$ws.Range("R4").Value = 7137099.645855149
$ws.Range("AW4").Value = "Jonathan Frendel"
$ws.Range("AX4").Value = "Jonathan Frendel, Astrid Blomberg, Elias Blad, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg"

# Row 5
$ws.Range("A5").Value = 111476580
$ws.Range("B5").Value = 89423
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P5").Value = "Sydvästra Gärdvattnet, Jmt"
$ws.Range("Q5").Value = 506446.2270308413
$ws.Range("R5").Value = 7137160.362918839
$ws.Range("AW5").Value = "Signe Propst"
$ws.Range("AX5").Value = "Signe Propst"

# Row 6
$ws.Range("A6").Value = 111477159
$ws.Range("B6").Value = 73696
$ws.Range("E6").Value = 6440
$ws.Range("F6").Value = "Vitgrynig nållav"
$ws.Range("G6").Value = "Chaenotheca subroscida"
$ws.Range("H6").Value = "(Eitner) Zahlbr."
$ws.Range("P6").Value = "SV gärdvattnet, Jmt"
$ws.Range("Q6").Value = 506548.1973901832
$ws.Range("R6").Value = 7137138.920001913
$ws.Range("AW6").Value = "Filippa Paperin"
$ws.Range("AX6").Value = "Filippa Paperin, Karl Soler Kinnerbäck, Tore Dahlberg, Melvin Lewin, Elvira Klang, Elicia Olsson, Jonathan Frendel, Astrid Blomberg, Iris Elmér, Ivar Anderberg, Kai Strömberg, Signe Propst, Elias Blad"

# Row 7
$ws.Range("A7").Value = 111479733
$ws.Range("B7").Value = 96265
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 219790
$ws.Range("F7").Value = "Fläcknycklar"
$ws.Range("G7").Value = "Dactylorhiza maculata"
$ws.Range("H7").Value = "(L.) Soó"
$ws.Range("P7").Value = "SV Gärdvattnet, Jmt"
$ws.Range("Q7").Value = 506542.631959103
$ws.Range("R7").Value = 7137104.68686779
$ws.Range("AW7").Value = "Jonathan Frendel"
$ws.Range("AX7").Value = "Jonathan Frendel, Astrid Blomberg, Elias Blad, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg"

# Row 8
$ws.Range("A8").Value = 111476591

# Row 9
$ws.Range("A9").Value = 111479730
$ws.Range("B9").Value = 77515
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 506636.7902023449
$ws.Range("R9").Value = 7137086.695334492

# Row 10
$ws.Range("A10").Value = 111479725
$ws.Range("B10").Value = 56398
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("M10").Value = "äldre spår"
$ws.Range("Q10").Value = 506573.8724125003
$ws.Range("R10").Value = 7137099.122253072

# Row 11
$ws.Range("A11").Value = 111476588
$ws.Range("B11").Value = 96368
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 221952
$ws.Range("F11").Value = "Spindelblomster"
$ws.Range("G11").Value = "Neottia cordata"
$ws.Range("H11").Value = "(L.) Rich."
$ws.Range("M11").Value = $null
$ws.Range("P11").Value = "Sydvästra Gärdvattnet, Jmt"
$ws.Range("Q11").Value = 506411.5719905405
$ws.Range("R11").Value = 7137139.931017525
$ws.Range("AW11").Value = "Signe Propst"
$ws.Range("AX11").Value = "Signe Propst"

# Row 12
$ws.Range("A12").Value = 111476577
$ws.Range("B12").Value = 96368
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 221952
$ws.Range("F12").Value = "Spindelblomster"
$ws.Range("G12").Value = "Neottia cordata"
$ws.Range("H12").Value = "(L.) Rich."
$ws.Range("P12").Value = "Sydvästra Gärdvattnet, Jmt"
$ws.Range("Q12").Value = 506399.4467280412
$ws.Range("R12").Value = 7137130.808956717
$ws.Range("AW12").Value = "Signe Propst"
$ws.Range("AX12").Value = "Signe Propst"

# Row 13
$ws.Range("A13").Value = 111479728
$ws.Range("B13").Value = 78579
$ws.Range("E13").Value = 2081
$ws.Range("F13").Value = "Skrovellav"
$ws.Range("G13").Value = "Lobaria scrobiculata"
$ws.Range("H13").Value = "(Scop.) DC."
$ws.Range("Q13").Value = 506536.1601844588
$ws.Range("R13").Value = 7137088.645264999

# Row 15
$ws.Range("A15").Value = 111479731
$ws.Range("B15").Value = 77515
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = "Garnlav"
$ws.Range("G15").Value = "Alectoria sarmentosa"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("P15").Value = "SV Gärdvattnet, Jmt"
$ws.Range("Q15").Value = 506570.8492999141
$ws.Range("R15").Value = 7137093.05119953
$ws.Range("AW15").Value = "Jonathan Frendel"
$ws.Range("AX15").Value = "Jonathan Frendel, Astrid Blomberg, Elias Blad, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg"

# Row 16
$ws.Range("A16").Value = 111479726
$ws.Range("B16").Value = 78605
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 6462
$ws.Range("F16").Value = "Stuplav"
$ws.Range("G16").Value = "Nephroma bellum"
$ws.Range("H16").Value = "(Spreng.) Tuck."
$ws.Range("M16").Value = $null
$ws.Range("Q16").Value = 506523.1528918216
$ws.Range("R16").Value = 7137086.451659708

# Row 17
$ws.Range("A17").Value = 111479729
$ws.Range("B17").Value = 78579
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 2081
$ws.Range("F17").Value = "Skrovellav"
$ws.Range("G17").Value = "Lobaria scrobiculata"
$ws.Range("H17").Value = "(Scop.) DC."
$ws.Range("Q17").Value = 506569.97720399
$ws.Range("R17").Value = 7137095.215254448

# Row 19
$ws.Range("A19").Value = 111480041
$ws.Range("B19").Value = 73696
$ws.Range("E19").Value = 6440
$ws.Range("F19").Value = "Vitgrynig nållav"
$ws.Range("G19").Value = "Chaenotheca subroscida"
$ws.Range("H19").Value = "(Eitner) Zahlbr."
$ws.Range("Q19").Value = 506532.6737455213
$ws.Range("R19").Value = 7137096.435164435

# Row 20
$ws.Range("A20").Value = 111480140
$ws.Range("B20").Value = 78605
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 6462
$ws.Range("F20").Value = "Stuplav"
$ws.Range("G20").Value = "Nephroma bellum"
$ws.Range("H20").Value = "(Spreng.) Tuck."
$ws.Range("Q20").Value = 506612.9201057266
$ws.Range("R20").Value = 7137093.574760968

# Row 21
$ws.Range("A21").Value = 111480182
$ws.Range("B21").Value = 78579
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 2081
$ws.Range("F21").Value = "Skrovellav"
$ws.Range("G21").Value = "Lobaria scrobiculata"
$ws.Range("H21").Value = "(Scop.) DC."
$ws.Range("Q21").Value = 506612.9247734078
$ws.Range("R21").Value = 7137091.40884747

# Row 22
$ws.Range("A22").Value = 111480035
$ws.Range("B22").Value = 76513
$ws.Range("E22").Value = 314
$ws.Range("F22").Value = "Vitskaftad svartspik"
$ws.Range("G22").Value = "Chaenothecopsis viridialba"
$ws.Range("H22").Value = "(Kremp.) A.F.W.Schmidt"
$ws.Range("Q22").Value = 506552.5373931379
$ws.Range("R22").Value = 7137137.629731925
